# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit diff:
#  - bumps the "Datos actualizados" timestamp
#  - updates per-country case/death figures for several rows
#  - re-ranks a few country pairs whose "Casos totales" crossed over,
#    which swaps both the country label and its row of figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 08:10"

# --- Row 27: Israel ---
$ws.Range("B27").Value = 166794
$ws.Range("C27").Value = 2392
$ws.Range("D27").Value = 122785
$ws.Range("E27").Value = 42862

# --- Rows 59/60: Uzbekistan overtakes Argelia in total cases, so they swap rank ---
# Row 59 becomes Uzbekistan with its freshly updated totals
$ws.Range("A59").Value = "Uzbekistan"
$ws.Range("B59").Value = 48776
$ws.Range("C59").Value = 347
$ws.Range("D59").Value = 45058
$ws.Range("E59").Value = 3313
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 405

# Row 60 becomes Argelia, keeping its previous (unchanged) totals
$ws.Range("A60").Value = "Argelia"
$ws.Range("B60").Value = 48734
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 34385
$ws.Range("E60").Value = 12717
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1632

# --- Row 64: Kirguistan ---
$ws.Range("B64").Value = 45072
$ws.Range("C64").Value = 73
$ws.Range("D64").Value = 41210
$ws.Range("E64").Value = 2799

# --- Row 75: El Salvador ---
$ws.Range("E75").Value = 6842
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 796

# --- Row 76: Australia ---
$ws.Range("B76").Value = 26779
$ws.Range("C76").Value = 40
$ws.Range("D76").Value = 23726
$ws.Range("E76").Value = 2229

# --- Rows 204/205: Timor Oriental / Santa Lucia swap places (tied totals, labels only) ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Rows 214/215: Islas Malvinas / Montserrat swap places (label + figures) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
